$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.222.04'
$ws.Range("E2").Value = '  +1.61%  '

$ws.Range("D3").Value = '2.419.23'
$ws.Range("E3").Value = '  +2.23%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '562.76'
$ws.Range("E5").Value = '  +2.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.27'
$ws.Range("E6").Value = '  +2.49%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E8").Value = '  +2.08%  '

$ws.Range("D9").Value = '2.416.00'
$ws.Range("E9").Value = '  +1.97%  '

$ws.Range("E10").Value = '  +1.69%  '

$ws.Range("E11").Value = '  -1.98%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.33'
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("E13").Value = '  +0.41%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.74'
$ws.Range("E14").Value = '  +0.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000175'
$ws.Range("E15").Value = '  +1.84%  '

$ws.Range("D16").Value = '2.851.91'
$ws.Range("E16").Value = '  +2.06%  '

$ws.Range("D17").Value = '62.148.13'
$ws.Range("E17").Value = '  +1.71%  '

$ws.Range("D18").Value = '2.414.28'
$ws.Range("E18").Value = '  +2.00%  '

$ws.Range("E19").Value = '  +2.37%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.19'
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '324.25'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.81'
$ws.Range("E22").Value = '  +2.50%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.71'
$ws.Range("E24").Value = '  +2.18%  '

$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.08'
$ws.Range("E26").Value = '  +2.28%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '578.93'
$ws.Range("E27").Value = '  +8.08%  '

$ws.Range("D28").Value = '0.0₃0953'
$ws.Range("E28").Value = '  +5.28%  '

$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.529.96'
$ws.Range("E29").Value = '  +2.27%  '

$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.23'
$ws.Range("E31").Value = '  -0.06%  '

$ws.Range("E32").Value = '  +3.11%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.149'
$ws.Range("E33").Value = '  +1.29%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.88'
$ws.Range("E34").Value = '  +1.93%  '

$ws.Range("E35").Value = '  +1.80%  '

$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '152.93'
$ws.Range("E39").Value = '  +4.90%  '

$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.382'
$ws.Range("E40").Value = '  +0.53%  '

$ws.Range("E41").Value = '  +1.02%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("E42").Value = '  -4.93%  '

$ws.Range("E43").Value = '  -0.44%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.30'
$ws.Range("E44").Value = '  +2.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.86'
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("E47").Value = '  +1.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.14'
$ws.Range("E48").Value = '  +0.16%  '

$ws.Range("E49").Value = '  +2.48%  '

$ws.Range("E50").Value = '  +1.52%  '

$ws.Range("E51").Value = '  +2.20%  '
